# Update the data source paths for the new repo structure
# (../mappings/data/... -> ../data/processed/...) and refresh the
# active-sheet / selection state left over from the editing session.

$wb = $excel.ActiveWorkbook

# --- Source sheet: rewrite the "source" file paths --------------------
$wsSource = $wb.Worksheets.Item("Source")

$papersOA = @("C2", "C5", "C8")
foreach ($addr in $papersOA) {
    $wsSource.Range($addr).Value = "../data/processed/papersOA.json"
}

$catalysts = @("C11", "C13", "C15", "C17", "C19", "C21", "C23", "C25", "C27", "C29", "C31", "C33", "C35", "C37", "C39", "C41", "C43", "C45", "C47", "C49", "C51", "C53")
foreach ($addr in $catalysts) {
    $wsSource.Range($addr).Value = "../data/processed/catalystsdata.csv"
}

$wsSource.Range("C55").Value = "../data/processed/chemicals.csv"

# --- Predicate_Object sheet: drop the stray fill style on a few cells -
$wsPO = $wb.Worksheets.Item("Predicate_Object")
$wsPO.Range("C72").Style = "Normal"
$wsPO.Range("C87").Style = "Normal"
$wsPO.Range("C88").Style = "Normal"

# --- Selection / active sheet bookkeeping ------------------------------
$wsSource.Activate()
$wsSource.Range("C2").Select()
